# Apply the two text edits described by the commit:
#  1) Slide 1 title: collapse the multi-run title into a single run with
#     new wording (keeps the paragraph's trailing endParaRPr formatting).
#  2) Slide 8 bullet: "Modul bluetooth yang digunakan adalah HC-05" ->
#     "Modul wifi yang digunakan adalah ESP8266".

$p = $ppt.ActivePresentation

# --- Slide 1: rewrite the title text -------------------------------------
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$titleLen = $titleRange.Length
$firstRunLen = 13   # "IMPLEMENTASI " (length of the run whose formatting we keep)

# Remove every run after the first one, leaving the trailing endParaRPr
# (the paragraph-mark formatting) untouched.
$remainder = $titleRange.Characters($firstRunLen + 1, $titleLen - $firstRunLen)
$remainder.Delete()

# Replace the remaining (first) run's text with the new title. The run
# keeps the original run's formatting (lang="en-US" sz=1800 b=1).
$titleShape.TextFrame.TextRange.Text = "PERMAINAN PENGENALAN KONSEP PEMROGRAMAN MENGGUNAKAN ROBOT KECERDASAN BUATAN PENELUSUR LABIRIN BERBASIS VISUAL PROGRAMMING DAN INTERNET OF THINGS"

# --- Slide 8: bluetooth/HC-05 -> wifi/ESP8266 -----------------------------
$s8 = $p.Slides.Item(8)
$specShape = $s8.Shapes.Item(2)
$specRange = $specShape.TextFrame2.TextRange

# Paragraph 3 of this shape reads "Modul bluetooth yang digunakan adalah HC-05"
$modulPara = $specRange.Paragraphs(3)
$btWord = $modulPara.Characters(7, 9)     # "bluetooth"
$btWord.Text = "wifi"

# Re-fetch the paragraph (its length/offsets shifted after the edit above)
# and replace the trailing " adalah HC-05" with " adalah ESP8266".
$modulPara2 = $specShape.TextFrame2.TextRange.Paragraphs(3)
$tail = $modulPara2.Characters(26, 13)    # " adalah HC-05"
$tail.Text = " adalah ESP8266"
